$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (matching the workbook's inline-string / shared-string
# storage) without leaving a stray number-format style attached to the cell, and without
# letting Excel auto-convert a numeric-looking string into a real number.
function Set-TextValue($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "247.61"
Set-TextValue $ws "D3" "21.97"
Set-TextValue $ws "D4" "5.353"
Set-TextValue $ws "D5" "0.05636"
Set-TextValue $ws "D6" "3.427"
Set-TextValue $ws "D7" "6.370"
Set-TextValue $ws "D8" "0.8185"
Set-TextValue $ws "D9" "0.9317"
Set-TextValue $ws "D10" "0.1438"
Set-TextValue $ws "D11" "0.07488"
Set-TextValue $ws "D13" "0.03081"
Set-TextValue $ws "D15" "3.565"
Set-TextValue $ws "D16" "0.001616"
Set-TextValue $ws "D18" "0.0005778"
Set-TextValue $ws "E18" "17OneONEWorstin24h"
Set-TextValue $ws "D20" "0.005063"
Set-TextValue $ws "D21" "0.001035"
Set-TextValue $ws "D22" "0.0001500"
Set-TextValue $ws "D24" "2.161"
Set-TextValue $ws "D25" "0.3306"
Set-TextValue $ws "D26" "0.1319"
Set-TextValue $ws "D28" "0.0002999"
Set-TextValue $ws "D40" "0.03954"
Set-TextValue $ws "D41" "0.006989"
Set-TextValue $ws "E41" "40KickTokenKICK"
Set-TextValue $ws "D42" "0.1065"
Set-TextValue $ws "D43" "0.003400"
Set-TextValue $ws "D44" "0.008541"
Set-TextValue $ws "D45" "0.00005578"
Set-TextValue $ws "D47" "0.0005498"
Set-TextValue $ws "D48" "0.7796"
Set-TextValue $ws "D49" "0.1770"
